$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking price cells so they keep their
# original textual representation (trailing zeros, thousands dots, etc.)
# instead of being coerced into floating point numbers.
$textForceCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D15", "D16", "D20", "D21", "D22", "D23", "D27", "D29", "D31", "D32", "D33", "D35", "D38", "D39", "D41", "D42", "D43", "D44", "D45", "D46", "D48", "D49")
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '43.597.94'
$ws.Range("E2").Value = '  -1.42%  '

$ws.Range("D3").Value = '2.290.69'
$ws.Range("E3").Value = '  +0.05%  '

$ws.Range("E4").Value = '  -0.20%  '

$ws.Range("D5").Value = '96.58'
$ws.Range("E5").Value = '  +0.49%  '

$ws.Range("D6").Value = '267.32'
$ws.Range("E6").Value = '  -2.61%  '

$ws.Range("D7").Value = '0.623'
$ws.Range("E7").Value = '  -0.79%  '

$ws.Range("E8").Value = '  -0.09%  '

$ws.Range("D9").Value = '0.610'
$ws.Range("E9").Value = '  -3.91%  '

$ws.Range("D10").Value = '45.32'
$ws.Range("E10").Value = '  -4.47%  '

$ws.Range("D11").Value = '0.0937'
$ws.Range("E11").Value = '  -0.33%  '

$ws.Range("E12").Value = '  -3.92%  '

$ws.Range("E13").Value = '  -0.10%  '

$ws.Range("D14").Value = '2.633.32'
$ws.Range("E14").Value = '  +0.16%  '

$ws.Range("D15").Value = '15.12'
$ws.Range("E15").Value = '  -3.50%  '

$ws.Range("D16").Value = '0.849'
$ws.Range("E16").Value = '  +1.07%  '

$ws.Range("D17").Value = '2.295.42'
$ws.Range("E17").Value = '  +0.13%  '

$ws.Range("D18").Value = '43.572.35'
$ws.Range("E18").Value = '  -1.57%  '

$ws.Range("E19").Value = '  +1.68%  '

$ws.Range("D20").Value = '6.17'
$ws.Range("E20").Value = '  -1.16%  '

$ws.Range("D21").Value = '71.93'
$ws.Range("E21").Value = '  +1.06%  '

$ws.Range("D22").Value = '2.41'
$ws.Range("E22").Value = '  +4.09%  '

$ws.Range("D23").Value = '232.73'
$ws.Range("E23").Value = '  -1.43%  '

$ws.Range("E24").Value = '  -10.48%  '

$ws.Range("E25").Value = '  -0.08%  '

$ws.Range("E26").Value = '  -1.23%  '

$ws.Range("D27").Value = '11.17'
$ws.Range("E27").Value = '  -3.05%  '

$ws.Range("E28").Value = '  +3.23%  '

$ws.Range("D29").Value = '40.46'
$ws.Range("E29").Value = '  +2.06%  '

$ws.Range("E30").Value = '  +0.89%  '

$ws.Range("D31").Value = '175.25'
$ws.Range("E31").Value = '  +0.91%  '

$ws.Range("D32").Value = '21.90'
$ws.Range("E32").Value = '  -0.56%  '

$ws.Range("D33").Value = '0.0884'
$ws.Range("E33").Value = '  -3.83%  '

$ws.Range("E34").Value = '  -5.34%  '

$ws.Range("D35").Value = '0.125'
$ws.Range("E35").Value = '  +0.40%  '

$ws.Range("E36").Value = '  -5.42%  '

$ws.Range("E37").Value = '  -0.83%  '

$ws.Range("D38").Value = '4.37'
$ws.Range("E38").Value = '  -1.83%  '

$ws.Range("D39").Value = '3.38'
$ws.Range("E39").Value = '  -4.26%  '

$ws.Range("E40").Value = '  -5.87%  '

$ws.Range("D41").Value = '2.32'
$ws.Range("E41").Value = '  +3.79%  '

$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").Value = '1.35'
$ws.Range("E42").Value = '  +13.97%  '

$ws.Range("B43").Value = 'Celestia'
$ws.Range("C43").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D43").Value = '12.22'
$ws.Range("E43").Value = '  -2.60%  '

$ws.Range("D44").Value = '63.64'
$ws.Range("E44").Value = '  +1.64%  '

$ws.Range("D45").Value = '8.79'
$ws.Range("E45").Value = '  +1.73%  '

$ws.Range("D46").Value = '5.24'
$ws.Range("E46").Value = '  -5.33%  '

$ws.Range("E47").Value = '  -0.71%  '

$ws.Range("D48").Value = '97.68'
$ws.Range("E48").Value = '  -2.96%  '

$ws.Range("D49").Value = '1.19'
$ws.Range("E49").Value = '  -0.55%  '

$ws.Range("D50").Value = '2.513.93'
$ws.Range("E50").Value = '  +0.18%  '

$ws.Range("E51").Value = '  -1.94%  '

# Restore the default cell style so we do not leave a stray "Text" number
# format applied to cells that originally had no explicit style.
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).Style = "Normal"
}
